# Add Azure IoT Hub device connection strings as a new header row on top of
# each device's data sheet, and update the saved selection on each sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "bin-008" --------------------------------------------------
$ws1 = $wb.Worksheets.Item("bin-008")
$ws1.Rows("1:1").Insert()
$ws1.Range("A1").Value = "HostName=filipiothub.azure-devices.net;DeviceId=bin-008;SharedAccessKey=Jm/zGBgY1ddtbRXaL9miMwi8zX/mHGMyBp1N8R4ZSYg="

# --- Sheet "bin-006" --------------------------------------------------
$ws2 = $wb.Worksheets.Item("bin-006")
$ws2.Rows("1:1").Insert()
$ws2.Range("A1").Value = "HostName=filipiothub.azure-devices.net;DeviceId=bin-006;SharedAccessKey=aNPACA6knFTNhaH5SN/WTsFh+Q6xF272WtJYef6RN2s="

# --- Update each sheet's saved selection/scroll position --------------
$ws2.Select()
$ws2.Range("E21").Select()

$ws1.Select()
$ws1.Range("G13").Select()
